# "string value formatting fix"
# Adds a new row (A14) on Sheet1 containing the literal text string
# "2/28/2024" (stored as a shared string, NOT converted to a date
# serial number), formatted with a date-looking number format, a thin
# box border around the cell, and left-aligned text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A14")

# Apply border + alignment first (order matters for how the
# number-format/style gets resolved on save).
$rng.HorizontalAlignment = -4131   # xlLeft
$rng.Borders.LineStyle = 1         # xlContinuous
$rng.Borders.Weight = 2            # xlThin

# Force the cell to Text format before writing the value so that the
# string "2/28/2024" is kept as literal text instead of being
# auto-parsed into a date serial number.
$rng.NumberFormat = "@"
$rng.Formula = "2/28/2024"

# Now apply the (date-style) display format on top of the text value.
$rng.NumberFormat = "mm-dd-yy"

[void]$rng.Select()
